# Generate Report for Archive
#
# The CI-generated localization-status report moved the single source
# document from "Ready for handoff" to "In Translation". That status
# string is shared by the Overview roll-up sheet (columns zh-cn/de-de)
# and by the per-locale detail sheets (zh-cn, de-de - "Status" column),
# so every occurrence of the old text needs to be updated in lock-step.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The Status columns were sized to fit "Ready for handoff"; with the
# shorter "In Translation" text the same columns re-fit narrower.
$overview.Columns.Item(5).ColumnWidth = 12.45
$overview.Columns.Item(6).ColumnWidth = 12.45
$zhcn.Columns.Item(3).ColumnWidth = 12.45
$dede.Columns.Item(3).ColumnWidth = 12.45
